$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '+87.94 ± 0.00'
$ws.Range("C2").Value = '+87.94 ± 0.00'
$ws.Range("D2").Value = '+6.15 ± 0.01'
$ws.Range("B3").Value = '+87.94 ± 0.00'
$ws.Range("C3").Value = '+87.94 ± 0.00'
$ws.Range("D4").Value = '+5.34 ± 0.01'
$ws.Range("D5").Value = '+0.81 ± 0.00'
$ws.Range("B6").Value = '-87.90 ± 0.01'
$ws.Range("C6").Value = '-87.87 ± 0.01'
$ws.Range("D6").Value = '-6.17 ± 0.01'
$ws.Range("C7").Value = '-5.34 ± 0.01'
$ws.Range("C8").Value = '-0.81 ± 0.00'
$ws.Range("B9").Value = '-3.32 ± 0.01'
$ws.Range("D9").Value = '-3.32 ± 0.01'
$ws.Range("B11").Value = '-80.68 ± 0.01'
$ws.Range("C11").Value = '-78.16 ± 0.01'
$ws.Range("D11").Value = '-2.52 ± 0.01'
$ws.Range("B12").Value = '-3.67 ± 0.00'
$ws.Range("C12").Value = '-3.56 ± 0.00'
$ws.Range("B13").Value = '+1.29 ± 0.00'
$ws.Range("C13").Value = '+1.19 ± 0.00'
$ws.Range("D13").Value = '+0.11 ± 0.00'
$ws.Range("B14").Value = '+5.37 ± 0.01'
$ws.Range("C14").Value = '+4.92 ± 0.01'
$ws.Range("D14").Value = '+6.37 ± 0.04'
